$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 146.41667
$ws.Range("I33").Value = 127.90909
$ws.Range("K33").Value = 127.90909
$ws.Range("M33").Value = 101.09091
$ws.Range("H80").Value = 946.5
$ws.Range("I80").Value = 838.8
$ws.Range("J80").Value = 1066.1666
$ws.Range("K80").Value = 2516.4
$ws.Range("L80").Value = 3198.4998
$ws.Range("M80").Value = -1518.4
$ws.Range("N80").Value = -5194.4998
$ws.Range("H83").Value = 946.5
$ws.Range("I83").Value = 838.8
$ws.Range("J83").Value = 1066.1666
$ws.Range("K83").Value = 7549.2
$ws.Range("L83").Value = 9595.499400000001
$ws.Range("M83").Value = -2557.2
$ws.Range("N83").Value = -19579.4994
$ws.Range("H98").Value = 1202.1154
$ws.Range("I98").Value = 875.2727
$ws.Range("J98").Value = 2999.75
$ws.Range("K98").Value = 875.2727
$ws.Range("L98").Value = 2999.75
$ws.Range("M98").Value = 622.7273
$ws.Range("N98").Value = -5995.75
$ws.Range("H112").Value = 1759
$ws.Range("J112").Value = 1593.5264
$ws.Range("L112").Value = 4780.5792
$ws.Range("N112").Value = -6996.5792
$ws.Range("H122").Value = 1202.1154
$ws.Range("I122").Value = 875.2727
$ws.Range("J122").Value = 2999.75
$ws.Range("K122").Value = 2625.8181
$ws.Range("L122").Value = 8999.25
$ws.Range("M122").Value = -175.8181
$ws.Range("N122").Value = -13899.25
$ws.Range("H137").Value = 8125.8696
$ws.Range("I137").Value = 1989.2858
$ws.Range("J137").Value = 17671.666
$ws.Range("K137").Value = 5967.857400000001
$ws.Range("L137").Value = 53014.99800000001
$ws.Range("M137").Value = -3417.857400000001
$ws.Range("N137").Value = -58114.99800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8742.896000000001
$ws.Range("I61").Value = 12141
$ws.Range("J61").Value = 6666.278
$ws.Range("K61").Value = 12141
$ws.Range("L61").Value = 6666.278
$ws.Range("M61").Value = -11929
$ws.Range("N61").Value = -7090.278
$ws.Range("H63").Value = 4443.625
$ws.Range("J63").Value = 7851.6665
$ws.Range("L63").Value = 7851.6665
$ws.Range("N63").Value = -9223.666499999999
$ws.Range("H66").Value = 4443.625
$ws.Range("J66").Value = 7851.6665
$ws.Range("L66").Value = 39258.3325
$ws.Range("N66").Value = -46122.3325
$ws.Range("H74").Value = 11357.839
$ws.Range("I74").Value = 9917.825999999999
$ws.Range("K74").Value = 9917.825999999999
$ws.Range("M74").Value = -9043.825999999999
$ws.Range("H77").Value = 11357.839
$ws.Range("I77").Value = 9917.825999999999
$ws.Range("K77").Value = 49589.13
$ws.Range("M77").Value = -45221.13
$ws.Range("H97").Value = 1379.2
$ws.Range("I97").Value = 1577.9
$ws.Range("J97").Value = 981.8
$ws.Range("K97").Value = 1577.9
$ws.Range("L97").Value = 981.8
$ws.Range("M97").Value = -1081.9
$ws.Range("N97").Value = -1973.8
$ws.Range("H110").Value = 1289.1351
$ws.Range("I110").Value = 974.9167
$ws.Range("K110").Value = 974.9167
$ws.Range("M110").Value = 1070.0833
$ws.Range("H132").Value = 5741.6787
$ws.Range("I132").Value = 4433.7856
$ws.Range("K132").Value = 13301.3568
$ws.Range("M132").Value = -10771.3568
$ws.Range("H136").Value = 8742.896000000001
$ws.Range("I136").Value = 12141
$ws.Range("J136").Value = 6666.278
$ws.Range("K136").Value = 36423
$ws.Range("L136").Value = 19998.834
$ws.Range("M136").Value = -33873
$ws.Range("N136").Value = -25098.834

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 113439.4
$ws.Range("I20").Value = 124812.664
$ws.Range("J20").Value = 11080
$ws.Range("K20").Value = 124812.664
$ws.Range("L20").Value = 11080
$ws.Range("M20").Value = -124565.664
$ws.Range("N20").Value = -11574
$ws.Range("H105").Value = 3068
$ws.Range("I105").Value = 1739.52
$ws.Range("J105").Value = 5835.6665
$ws.Range("K105").Value = 1739.52
$ws.Range("L105").Value = 5835.6665
$ws.Range("M105").Value = 7.480000000000018
$ws.Range("N105").Value = -9329.666499999999
$ws.Range("H134").Value = 6255
$ws.Range("I134").Value = 6862.9287
$ws.Range("K134").Value = 20588.7861
$ws.Range("M134").Value = -18053.7861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 71794.92999999999
$ws.Range("I16").Value = 47894.25
$ws.Range("J16").Value = 167397.67
$ws.Range("K16").Value = 47894.25
$ws.Range("L16").Value = 167397.67
$ws.Range("M16").Value = -47607.25
$ws.Range("N16").Value = -167971.67
$ws.Range("H58").Value = 4533.0967
$ws.Range("I58").Value = 2744.5386
$ws.Range("K58").Value = 2744.5386
$ws.Range("M58").Value = -2541.5386
$ws.Range("H59").Value = 33399
$ws.Range("I59").Value = 32497.5
$ws.Range("J59").Value = 34000
$ws.Range("K59").Value = 32497.5
$ws.Range("L59").Value = 34000
$ws.Range("M59").Value = -31352.5
$ws.Range("N59").Value = -36290
$ws.Range("H68").Value = 40827.223
$ws.Range("J68").Value = 40827.223
$ws.Range("L68").Value = 40827.223
$ws.Range("N68").Value = -42325.223
$ws.Range("H71").Value = 40827.223
$ws.Range("J71").Value = 40827.223
$ws.Range("L71").Value = 122481.669
$ws.Range("N71").Value = -129969.669
$ws.Range("H74").Value = 37889.25
$ws.Range("J74").Value = 37889.25
$ws.Range("L74").Value = 37889.25
$ws.Range("N74").Value = -39637.25
$ws.Range("H77").Value = 37889.25
$ws.Range("J77").Value = 37889.25
$ws.Range("L77").Value = 113667.75
$ws.Range("N77").Value = -122403.75
$ws.Range("H113").Value = 71794.92999999999
$ws.Range("I113").Value = 47894.25
$ws.Range("J113").Value = 167397.67
$ws.Range("K113").Value = 47894.25
$ws.Range("L113").Value = 167397.67
$ws.Range("M113").Value = -45724.25
$ws.Range("N113").Value = -171737.67
$ws.Range("H132").Value = 10257.435
$ws.Range("I132").Value = 10257.435
$ws.Range("K132").Value = 30772.305
$ws.Range("M132").Value = -28242.305
$ws.Range("H134").Value = 3857.862
$ws.Range("I134").Value = 3584.3704
$ws.Range("K134").Value = 10753.1112
$ws.Range("M134").Value = -8218.111199999999
$ws.Range("H136").Value = 4533.0967
$ws.Range("I136").Value = 2744.5386
$ws.Range("K136").Value = 8233.6158
$ws.Range("M136").Value = -5683.6158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 1116.7273
$ws.Range("I50").Value = 322.85715
$ws.Range("K50").Value = 968.5714499999999
$ws.Range("M50").Value = -487.5714499999999
$ws.Range("H53").Value = 1116.7273
$ws.Range("I53").Value = 322.85715
$ws.Range("K53").Value = 968.5714499999999
$ws.Range("M53").Value = -487.5714499999999
$ws.Range("H116").Value = 5236
$ws.Range("I116").Value = 2226.6667
$ws.Range("K116").Value = 6680.000100000001
$ws.Range("M116").Value = -3238.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H102").Value = 1043.75
$ws.Range("I102").Value = 1043.75
$ws.Range("K102").Value = 1043.75
$ws.Range("M102").Value = 578.25
$ws.Range("H113").Value = 1500
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 3344
$ws.Range("I122").Value = 1250
$ws.Range("J122").Value = 4181.6
$ws.Range("K122").Value = 3750
$ws.Range("L122").Value = 12544.8
$ws.Range("M122").Value = -1300
$ws.Range("N122").Value = -17444.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4834.9614
$ws.Range("I132").Value = 3304.2144
$ws.Range("J132").Value = 6620.8335
$ws.Range("K132").Value = 9912.643199999999
$ws.Range("L132").Value = 19862.5005
$ws.Range("M132").Value = -7382.643199999999
$ws.Range("N132").Value = -24922.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7418.381
$ws.Range("I81").Value = 1941.1666
$ws.Range("K81").Value = 3882.3332
$ws.Range("M81").Value = -2821.3332
$ws.Range("H84").Value = 7418.381
$ws.Range("I84").Value = 1941.1666
$ws.Range("K84").Value = 19411.666
$ws.Range("M84").Value = -14107.666
$ws.Range("H122").Value = 50313.74
$ws.Range("I122").Value = 1167.8422
$ws.Range("K122").Value = 3503.5266
$ws.Range("M122").Value = -1053.5266
$ws.Range("H126").Value = 1662.2354
$ws.Range("I126").Value = 1392.9231
$ws.Range("K126").Value = 4178.7693
$ws.Range("M126").Value = -1708.7693
$ws.Range("H136").Value = 52486.094
$ws.Range("I136").Value = 5110.4
$ws.Range("K136").Value = 15331.2
$ws.Range("M136").Value = -12781.2
